# Mise à jour de extimation
# Apply the same edits that were made in the tracked commit:
#  - On sheet "Iteration #2", fill in rows 17 and 18 of the work log with
#    new dates / tasks / hours (previously blank placeholder rows).
#  - The shared-string table, the SUM total, and the view's
#    selection/scroll position are all recomputed automatically by Excel
#    when the cell contents change.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Iteration #2")

# Row 17: new task entry (replaces the previously blank placeholder row
# that only had a stray space character in column A). Use the raw Excel
# date serial number (instead of a DateTime) so the existing date style
# from A16 is reused as-is, without Excel inventing a new number format.
$ws.Range("A17").Value = 43202
$ws.Range("A17").NumberFormat = $ws.Range("A16").NumberFormat
$ws.Range("B17").Value = "Debug service + fin debug asynctask"
$ws.Range("C17").Value = 4

# Row 18: new task entry (previously completely empty)
$ws.Range("A18").Value = 43206
$ws.Range("A18").NumberFormat = $ws.Range("A16").NumberFormat
$ws.Range("B18").Value = "fin debug service + aide francois-oli"
$ws.Range("C18").Value = 3

# Match the view state left behind after the edit (scrolled down a bit,
# selection left on B20)
$ws.Activate()
$ws.Application.ActiveWindow.ScrollRow = 11
$ws.Range("B20").Select()
